$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-9, columns B (message), C (sentence), D (boundary)
$data = @(
    @{Row=2; B="Compass Error Compass Error Compass disconnected ."; C="Compass Error"; D="0-1"},
    @{Row=3; B="Compass Error Compass Error Compass disconnected ."; C="Compass Error"; D="2-3"},
    @{Row=4; B="Compass Error Compass Error Compass disconnected ."; C="Compass disconnected"; D="4-5"},
    @{Row=5; B="Compass Error Compass Error Compass disconnected ."; C="Compass Error Compass Error Compass disconnected"; D="0-5"},
    @{Row=6; B="Critical low battery Aircraft in Auto Power Off Protection Forced landing in progress ."; C="Critical low battery"; D="0-2"},
    @{Row=7; B="Critical low battery Aircraft in Auto Power Off Protection Forced landing in progress ."; C="Aircraft in Auto Power Off Protection"; D="3-8"},
    @{Row=8; B="Critical low battery Aircraft in Auto Power Off Protection Forced landing in progress ."; C="Forced landing in progress"; D="9-12"},
    @{Row=9; B="Critical low battery Aircraft in Auto Power Off Protection Forced landing in progress ."; C="Critical low battery Aircraft in Auto Power Off Protection Forced landing in progress"; D="0-12"}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
}
